# Update the footer "Date and time" placeholder text from 2020/6/20 to
# 2020/6/21 everywhere it appears: every slide, every slide layout, the
# slide master, and the notes master.

$p = $ppt.ActivePresentation
$oldDate = "2020/6/20"
$newDate = "2020/6/21"

function Update-DateAndTime($headersFooters) {
    $dt = $headersFooters.DateAndTime
    if ($dt -ne $null) {
        $dt.Text = $newDate
    }
}

# 1) Every slide's footer date field.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    Update-DateAndTime $slide.HeadersFooters
}

# 2) Slide master's footer date field.
$master = $p.SlideMaster
Update-DateAndTime $master.HeadersFooters

# 3) Every slide layout (custom layout) attached to the slide master.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DateAndTime $layout.HeadersFooters
}

# 4) Notes master's footer date field.
$notesMaster = $p.NotesMaster
Update-DateAndTime $notesMaster.HeadersFooters
